$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New departure rows to append to the "Main Data" sheet
$rows = @(
    @("74", "Saturday, Jan 14", "9:50 PM",  "W95154", "London",          "(LTN)", "Wizz Air ",                                     "A320", "(G-WUKF)", "10:07 PM", "0 hours, 17 minutes"),
    @("75", "Saturday, Jan 14", "10:05 PM", "FR3473", "London",          "(LTN)", "Ryanair ",                                      "B738", "(EI-EFJ)", "10:31 PM", "0 hours, 26 minutes"),
    @("76", "Sunday, Jan 15",   "12:05 AM", "PQ7552", "Sharm el-Sheikh", "(SSH)", "SkyUp Airlines (The Power Of Freedom Livery) ", "B738", "(UR-SQM)", "12:27 AM", "0 hours, 22 minutes")
)

$startRow = 75
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = [double]$data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    # Column K (11) is left blank, matching the existing rows above it
    $ws.Cells.Item($r, 11).Borders.LineStyle = -4142
    $ws.Cells.Item($r, 12).Value = $data[10]
    # Column M (13) is left blank, matching the existing rows above it
    $ws.Cells.Item($r, 13).Borders.LineStyle = -4142
}
